$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 810.3111
$ws.Range("I15").Value = 810.3111
$ws.Range("K15").Value = 2430.9333
$ws.Range("M15").Value = -2261.9333
$ws.Range("H76").Value = 5642
$ws.Range("I76").Value = 5298.8
$ws.Range("K76").Value = 5298.8
$ws.Range("M76").Value = -4983.8
$ws.Range("H79").Value = 5642
$ws.Range("I79").Value = 5298.8
$ws.Range("K79").Value = 5298.8
$ws.Range("M79").Value = -4206.8
$ws.Range("H86").Value = 5973.92
$ws.Range("I86").Value = 3920.4
$ws.Range("J86").Value = 9054.200000000001
$ws.Range("K86").Value = 3920.4
$ws.Range("L86").Value = 9054.200000000001
$ws.Range("M86").Value = -2797.4
$ws.Range("N86").Value = -11300.2
$ws.Range("H89").Value = 5973.92
$ws.Range("I89").Value = 3920.4
$ws.Range("J89").Value = 9054.200000000001
$ws.Range("K89").Value = 19602
$ws.Range("L89").Value = 45271
$ws.Range("M89").Value = -13986
$ws.Range("N89").Value = -56503
$ws.Range("H106").Value = 5442.2144
$ws.Range("J106").Value = 5616.3335
$ws.Range("L106").Value = 5616.3335
$ws.Range("N106").Value = -6878.3335

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 310.83334
$ws.Range("I4").Value = 273
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 273
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -157
$ws.Range("N4").Value = -732
$ws.Range("H25").Value = 6611
$ws.Range("I25").Value = 1500
$ws.Range("J25").Value = 9166.5
$ws.Range("K25").Value = 1500
$ws.Range("L25").Value = 9166.5
$ws.Range("M25").Value = -1098
$ws.Range("N25").Value = -9970.5
$ws.Range("H32").Value = 8261.25
$ws.Range("I32").Value = 5293.8335
$ws.Range("J32").Value = 23840.188
$ws.Range("K32").Value = 5293.8335
$ws.Range("L32").Value = 23840.188
$ws.Range("M32").Value = -5006.8335
$ws.Range("N32").Value = -24414.188
$ws.Range("H43").Value = 25377
$ws.Range("J43").Value = 25377
$ws.Range("L43").Value = 25377
$ws.Range("N43").Value = -26003
$ws.Range("H45").Value = 9779.177
$ws.Range("I45").Value = 11322.454
$ws.Range("J45").Value = 6949.8335
$ws.Range("K45").Value = 11322.454
$ws.Range("L45").Value = 6949.8335
$ws.Range("M45").Value = -10945.454
$ws.Range("N45").Value = -7703.8335
$ws.Range("H61").Value = 2231.9
$ws.Range("I61").Value = 2007.6666
$ws.Range("J61").Value = 4250
$ws.Range("K61").Value = 2007.6666
$ws.Range("L61").Value = 4250
$ws.Range("M61").Value = -1795.6666
$ws.Range("N61").Value = -4674
$ws.Range("H63").Value = 2812
$ws.Range("I63").Value = 1968.25
$ws.Range("K63").Value = 1968.25
$ws.Range("M63").Value = -1282.25
$ws.Range("H66").Value = 2812
$ws.Range("I66").Value = 1968.25
$ws.Range("K66").Value = 9841.25
$ws.Range("M66").Value = -6409.25
$ws.Range("H74").Value = 2832.923
$ws.Range("I74").Value = 2575.2727
$ws.Range("J74").Value = 4250
$ws.Range("K74").Value = 2575.2727
$ws.Range("L74").Value = 4250
$ws.Range("M74").Value = -1701.2727
$ws.Range("N74").Value = -5998
$ws.Range("H77").Value = 2832.923
$ws.Range("I77").Value = 2575.2727
$ws.Range("J77").Value = 4250
$ws.Range("K77").Value = 12876.3635
$ws.Range("L77").Value = 21250
$ws.Range("M77").Value = -8508.363499999999
$ws.Range("N77").Value = -29986
$ws.Range("H110").Value = 1379.909
$ws.Range("J110").Value = 1807.25
$ws.Range("L110").Value = 1807.25
$ws.Range("N110").Value = -5897.25
$ws.Range("H122").Value = 5261.875
$ws.Range("I122").Value = 5335.6875
$ws.Range("J122").Value = 5114.25
$ws.Range("K122").Value = 16007.0625
$ws.Range("L122").Value = 15342.75
$ws.Range("M122").Value = -13557.0625
$ws.Range("N122").Value = -20242.75
$ws.Range("H132").Value = 1812.6666
$ws.Range("I132").Value = 1746.826
$ws.Range("K132").Value = 5240.478
$ws.Range("M132").Value = -2710.478
$ws.Range("H136").Value = 2231.9
$ws.Range("I136").Value = 2007.6666
$ws.Range("J136").Value = 4250
$ws.Range("K136").Value = 6022.9998
$ws.Range("L136").Value = 12750
$ws.Range("M136").Value = -3472.9998
$ws.Range("N136").Value = -17850

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2234922.5
$ws.Range("I99").Value = 2917.6
$ws.Range("K99").Value = 2917.6
$ws.Range("M99").Value = -1419.6
$ws.Range("H112").Value = 62540
$ws.Range("J112").Value = 67548.89
$ws.Range("L112").Value = 67548.89
$ws.Range("N112").Value = -70502.89
$ws.Range("H134").Value = 2624.75
$ws.Range("J134").Value = 3083.3333
$ws.Range("L134").Value = 9249.999899999999
$ws.Range("N134").Value = -14319.9999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5181
$ws.Range("I31").Value = 3174.4285
$ws.Range("J31").Value = 6087.1934
$ws.Range("K31").Value = 3174.4285
$ws.Range("L31").Value = 6087.1934
$ws.Range("M31").Value = -2879.4285
$ws.Range("N31").Value = -6677.1934
$ws.Range("H34").Value = 5181
$ws.Range("I34").Value = 3174.4285
$ws.Range("J34").Value = 6087.1934
$ws.Range("K34").Value = 3174.4285
$ws.Range("L34").Value = 6087.1934
$ws.Range("M34").Value = -2972.4285
$ws.Range("N34").Value = -6491.1934
$ws.Range("H107").Value = 1048.8
$ws.Range("I107").Value = 943.4
$ws.Range("K107").Value = 943.4
$ws.Range("M107").Value = 976.6
$ws.Range("H122").Value = 4240.6665
$ws.Range("I122").Value = 3512.5715
$ws.Range("K122").Value = 10537.7145
$ws.Range("M122").Value = -8087.7145
$ws.Range("H131").Value = 38331.5
$ws.Range("J131").Value = 38331.5
$ws.Range("L131").Value = 38331.5
$ws.Range("N131").Value = -48411.5
$ws.Range("H132").Value = 2217.0588
$ws.Range("I132").Value = 1884.2727
$ws.Range("K132").Value = 5652.8181
$ws.Range("M132").Value = -3122.8181
$ws.Range("H134").Value = 1813.1666
$ws.Range("I134").Value = 1506.24
$ws.Range("K134").Value = 4518.72
$ws.Range("M134").Value = -1983.72

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1430.7931
$ws.Range("I102").Value = 1517.84
$ws.Range("K102").Value = 1517.84
$ws.Range("M102").Value = 104.1600000000001
$ws.Range("H132").Value = 1871.9565
$ws.Range("I132").Value = 1540.25
$ws.Range("K132").Value = 4620.75
$ws.Range("M132").Value = -2090.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1386.24
$ws.Range("I16").Value = 1533.2727
$ws.Range("J16").Value = 308
$ws.Range("K16").Value = 1533.2727
$ws.Range("L16").Value = 308
$ws.Range("M16").Value = -1363.2727
$ws.Range("N16").Value = -648
$ws.Range("H46").Value = 1788.037
$ws.Range("I46").Value = 1550.6154
$ws.Range("J46").Value = 2008.5
$ws.Range("K46").Value = 1550.6154
$ws.Range("L46").Value = 2008.5
$ws.Range("M46").Value = -1362.6154
$ws.Range("N46").Value = -2384.5
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 8009154.5
$ws.Range("I122").Value = 8692.923000000001
$ws.Range("K122").Value = 26078.769
$ws.Range("M122").Value = -23628.769
$ws.Range("H131").Value = 52995.668
$ws.Range("I131").Value = 50000
$ws.Range("J131").Value = 54493.5
$ws.Range("K131").Value = 50000
$ws.Range("L131").Value = 54493.5
$ws.Range("M131").Value = -44960
$ws.Range("N131").Value = -64573.5
$ws.Range("H132").Value = 2788.3572
$ws.Range("I132").Value = 2753.3333
$ws.Range("K132").Value = 8259.999899999999
$ws.Range("M132").Value = -5729.999899999999
$ws.Range("H136").Value = 9656.235000000001
$ws.Range("I136").Value = 10515.5
$ws.Range("K136").Value = 31546.5
$ws.Range("M136").Value = -28996.5
$ws.Range("H139").Value = 79606.42999999999
$ws.Range("J139").Value = 89561.25
$ws.Range("L139").Value = 89561.25
$ws.Range("N139").Value = -99841.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 26990
$ws.Range("J26").Value = 26990
$ws.Range("L26").Value = 26990
$ws.Range("N26").Value = -27576
$ws.Range("H132").Value = 1404451.8
$ws.Range("I132").Value = 1428.7084
$ws.Range("K132").Value = 4286.1252
$ws.Range("M132").Value = -1756.1252
$ws.Range("H139").Value = 92334.27
$ws.Range("J139").Value = 92334.27
$ws.Range("L139").Value = 92334.27
$ws.Range("N139").Value = -102614.27
